$d = $word.ActiveDocument

# Locate the paragraph "One that takes a mirror for calls to ObjectMirage.make"
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    if ($pp.Range.Text.StartsWith("One that takes a mirror for calls to ObjectMirage.make")) {
        $target = $pp
        break
    }
}

$nextPara = $target.Next()
$insertStart = $nextPara.Range.Start
$insertPoint = $d.Range($insertStart, $insertStart)

# Capture the exact text of the paragraph that currently follows (so we can
# rebuild it verbatim as the trailing element of the spliced-in OOXML; this
# is required because InsertXML always folds the *last* <w:p> of the inserted
# fragment into the paragraph located at the insertion point).
$followingText = $nextPara.Range.Text
# Range.Text includes the trailing paragraph mark; strip it.
$followingText = $followingText.TrimEnd([char]13, [char]7)

$xml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="3"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Problem </w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="4"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>In Java, a</w:t></w:r><w:r><w:t>nonymous inner class constructors will set fields before calling super</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="4"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Bytecode in general can do this whenever it wants</w:t></w:r></w:p>
<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">$followingText</w:t></w:r></w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insertPoint.InsertXML($xml)

# InsertXML merges the final <w:p> above with the paragraph that used to sit
# at $insertStart, duplicating its text ("Static fields" + "Static fields").
# Remove the now-redundant first copy, leaving the original paragraph intact.
$dupPara = $d.Paragraphs.Item($insertPoint)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs.Item($i)
    $txt = $pp.Range.Text.TrimEnd([char]13, [char]7)
    if ($txt -eq ($followingText + $followingText)) {
        $dupPara = $pp
        break
    }
}
$delRange = $d.Range($dupPara.Range.Start, $dupPara.Range.Start + $followingText.Length)
$delRange.Delete()
